$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 10; existing rows 10..106 shift down to 11..107
$ws.Rows.Item(10).Insert()

# Populate the newly inserted row 10 with the new data record
$ws.Range("A10").Value = 3
$ws.Range("B10").Value = "Femacal de La Calera"
$ws.Range("C10").Value = "Coquimbo"
$ws.Range("D10").Value = 45163
$ws.Range("E10").Value = 5
$ws.Range("F10").Value = 100112035
$ws.Range("G10").Value = "Bruselas (repollito)"
$ws.Range("H10").Value = "Sin especificar"
$ws.Range("I10").Value = "Primera"
$ws.Range("J10").Value = 40
$ws.Range("K10").Value = 11000
$ws.Range("L10").Value = 11000
$ws.Range("M10").Value = 11000
$ws.Range("N10").Value = "$/malla 10 kilos"
$ws.Range("O10").Value = "Provincia de Quillota"
$ws.Range("P10").Value = 1100
$ws.Range("Q10").Value = 10
$ws.Range("R10").Value = "Hortaliza"
